$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.269.09"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "3.811.97"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "702.33"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.98"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("D7").Value = "3.811.93"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.47"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.487"
$ws.Range("E12").Value = "  +5.98%  "
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.67"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "4.457.72"
$ws.Range("D16").Value = "3.798.99"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "71.382.40"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.23"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.52"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "511.67"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.46"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.715"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.81"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.76"
$ws.Range("E26").Value = "  +4.70%  "
$ws.Range("D27").Value = "3.959.16"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.32"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -4.60%  "
$ws.Range("E31").Value = "  -5.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.27"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.35"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.21"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("E35").Value = "  -5.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.32"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").Value = "3.776.54"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.65"
$ws.Range("E39").Value = "  +10.59%  "
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("E41").Value = "  +5.21%  "
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.22"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "166.08"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "50.06"
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "432.45"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("E49").Value = "  -5.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.92"
$ws.Range("E50").Value = "  +9.78%  "
$ws.Range("E51").Value = "  -0.40%  "
